# CryCompanywiseStockReport_1 - correct mis-paired stock rows.
#
# The source report pairs up two rows per SKU (one line per incoming batch /
# rate); a transcription bug swapped the Item Code / MRP / Qty / Value
# (columns B, C, D, E, F, G) between adjacent rows for a number of SKUs
# across the sheet. This script restores the correct pairing by swapping
# columns B:G between the affected rows (A - the serial number - and H:M -
# all blank - are left untouched).
#
# One group (rows 408-410) is a three-way rotation rather than a simple
# swap, handled separately below.
#
# NOTE: this engine's PowerShell function calls only bind parameters
# positionally (named args / default param values are not honoured), so
# all helper calls below pass every argument explicitly, in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowColumns {
    param($Row1, $Row2, $FirstCol, $LastCol)

    for ($col = $FirstCol; $col -le $LastCol; $col++) {
        $cell1 = $ws.Cells.Item($Row1, $col)
        $cell2 = $ws.Cells.Item($Row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Simple pairwise swaps (adjacent rows whose B:G data got transposed).
$swapPairs = @(
    @(49, 50),
    @(369, 370),
    @(427, 428),
    @(435, 436),
    @(438, 439),
    @(468, 469),
    @(497, 498),
    @(602, 603),
    @(616, 617),
    @(620, 621),
    @(782, 783),
    @(831, 832),
    @(841, 842),
    @(843, 844),
    @(845, 846),
    @(870, 871),
    @(872, 873),
    @(878, 879),
    @(887, 888),
    @(896, 897),
    @(902, 903),
    @(904, 905),
    @(939, 940),
    @(946, 947),
    @(977, 978)
)

foreach ($pair in $swapPairs) {
    Swap-RowColumns $pair[0] $pair[1] 2 7
}

# Three-way rotation: new row408 <- old row409, new row409 <- old row410,
# new row410 <- old row408.
$firstCol = 2  # B
$lastCol = 7   # G
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $c408 = $ws.Cells.Item(408, $col)
    $c409 = $ws.Cells.Item(409, $col)
    $c410 = $ws.Cells.Item(410, $col)

    $v408 = $c408.Value2
    $v409 = $c409.Value2
    $v410 = $c410.Value2

    $c408.Value = $v409
    $c409.Value = $v410
    $c410.Value = $v408
}
